$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data (filtered save games) for rows 2-8, columns B,C,D,E,G
# Column F (Win) is left unchanged.

$data = @{
    2 = @{ B = 0.04763786555579896; C = 0.3127903958511391;  D = 0.8054896365839992;  E = 8.660232485948974; G = 9.826150383939911 }
    3 = @{ B = 0.6753301551942219;  C = 1.667794583268128;   D = 0.1575252929769615;  E = 0.496779210170732; G = 2.997429241610044 }
    4 = @{ B = 0.3048080303191223;  C = 10.29869402782916;   D = 0.1575252929769615;  E = 8.660232485948974; G = 19.42125983707422 }
    5 = @{ B = 0.6753301551942219;  C = 1.667794583268128;   D = 0.1575252929769615;  E = 8.660232485948974; G = 11.16088251738829 }
    6 = @{ B = 3.230985683306322;   C = 1.667794583268128;   D = 0.8054896365839992;  E = 0.496779210170732; G = 6.201049113329182 }
    7 = @{ B = 3.230985683306322;   C = 1.667794583268128;   D = 3.900430680208489;   E = 8.660232485948974; G = 17.45944343273191 }
    8 = @{ B = 1.459612070389937;   C = 1.667794583268128;   D = 0.8054896365839992;  E = 0.496779210170732; G = 4.429675500412797 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
